$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2022-06-25"

# Update the month/"through" label used in column B's header
$ws.Range("B1").Value = "June 2022 (through June 25)"

# Update existing counts that changed (new carjacking incidents for 2022-07-03)
$ws.Range("B5").Value = 3    # South Shore, June 2022 (through June 25)
$ws.Range("H6").Value = 5    # Humboldt Park, June 2021
$ws.Range("B10").Value = 7   # Garfield Park, June 2022 (through June 25)
$ws.Range("H12").Value = 2   # Roseland, June 2021
$ws.Range("B20").Value = 3   # Hyde Park, June 2022 (through June 25)

# Fill in previously-empty cells with new counts
$ws.Range("Z27").Value = 1   # United Center, June 2018
$ws.Range("H28").Value = 1   # Ukrainian Village, June 2021
$ws.Range("H41").Value = 1   # Woodlawn, June 2021
$ws.Range("N45").Value = 1   # Avalon Park, June 2020
$ws.Range("AF46").Value = 1  # Avondale, June 2017
$ws.Range("H55").Value = 1   # East Village, June 2021
$ws.Range("H80").Value = 1   # Oakland, June 2021
$ws.Range("H95").Value = 1   # West Town, June 2021
